# Weekly update: insert 4 new price rows at the top of this week's block
# (rows 325-328), pushing the previously-last week's rows (325-338) down
# to (329-342).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows before row 325 (existing rows 325-338 shift to 329-342).
$ws.Range("A325:A328").EntireRow.Insert()

# New row 325
$ws.Range("A325").Value = 4
$ws.Range("B325").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C325").Value = "Los Lagos"
$ws.Range("D325").Value = 44509
$ws.Range("E325").Value = 10
$ws.Range("F325").Value = 100112004
$ws.Range("G325").Value = "Cebolla"
$ws.Range("H325").Value = "Morada(o)"
$ws.Range("I325").Value = "1a nueva(o)"
$ws.Range("J325").Value = 220
$ws.Range("K325").Value = 11000
$ws.Range("L325").Value = 11000
$ws.Range("M325").Value = 11000
$ws.Range("N325").Value = "`$/malla 18 kilos"
$ws.Range("O325").Value = "Región de O'Higgins"
$ws.Range("P325").Value = 611
$ws.Range("Q325").Value = 18
$ws.Range("R325").Value = "Hortaliza"

# New row 326
$ws.Range("A326").Value = 4
$ws.Range("B326").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C326").Value = "Los Lagos"
$ws.Range("D326").Value = 44509
$ws.Range("E326").Value = 10
$ws.Range("F326").Value = 100112004
$ws.Range("G326").Value = "Cebolla"
$ws.Range("H326").Value = "Sin especificar"
$ws.Range("I326").Value = "1a nueva(o)"
$ws.Range("J326").Value = 900
$ws.Range("K326").Value = 6500
$ws.Range("L326").Value = 7000
$ws.Range("M326").Value = 6667
$ws.Range("N326").Value = "`$/malla 16 kilos"
$ws.Range("O326").Value = "Región de O'Higgins"
$ws.Range("P326").Value = 417
$ws.Range("Q326").Value = 16
$ws.Range("R326").Value = "Hortaliza"

# New row 327
$ws.Range("A327").Value = 4
$ws.Range("B327").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C327").Value = "Los Lagos"
$ws.Range("D327").Value = 44509
$ws.Range("E327").Value = 10
$ws.Range("F327").Value = 100112004
$ws.Range("G327").Value = "Cebolla"
$ws.Range("H327").Value = "Sin especificar"
$ws.Range("I327").Value = "1a nueva(o)"
$ws.Range("J327").Value = 300
$ws.Range("K327").Value = 9000
$ws.Range("L327").Value = 9000
$ws.Range("M327").Value = 9000
$ws.Range("N327").Value = "`$/malla 18 kilos"
$ws.Range("O327").Value = "Región de O'Higgins"
$ws.Range("P327").Value = 500
$ws.Range("Q327").Value = 18
$ws.Range("R327").Value = "Hortaliza"

# New row 328
$ws.Range("A328").Value = 4
$ws.Range("B328").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C328").Value = "Los Lagos"
$ws.Range("D328").Value = 44509
$ws.Range("E328").Value = 10
$ws.Range("F328").Value = 100112004
$ws.Range("G328").Value = "Cebolla"
$ws.Range("H328").Value = "Sin especificar"
$ws.Range("I328").Value = "Primera"
$ws.Range("J328").Value = 300
$ws.Range("K328").Value = 9000
$ws.Range("L328").Value = 9000
$ws.Range("M328").Value = 9000
$ws.Range("N328").Value = "`$/malla 18 kilos"
$ws.Range("O328").Value = "Perú"
$ws.Range("P328").Value = 500
$ws.Range("Q328").Value = 18
$ws.Range("R328").Value = "Hortaliza"
